$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)
$shp = $s.Shapes.Item("Rectangle 7")
# Narrow the "TweetContent (as a form)" box: cx 8264434 -> 8122891 EMU
# (1 point = 12700 EMU, so 8122891 / 12700 = 639.597716535433 pt)
$shp.Width = 639.597716535433
